# Automatische test-sync: 2025-07-29 22:09:50
# Adds Testmail #19 ("Bel jij klant Jansen even?") as a new row to the
# "Logs" sheet and updates the "Dashboard" summary sheet + conditional
# formatting ranges to account for it.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 21 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A21").Value = "Bel jij klant Jansen even?"
$logs.Range("B21").Value = "mailmind.test@zohomail.eu"
$logs.Range("C21").Value = "Testmail #19: Bel jij klant Jansen even?"
$logs.Range("D21").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("F21").Value = "2025-07-29 22:08:55"
$logs.Range("G21").Value = "Nee"
$logs.Range("H21").Value = "Ja"
$logs.Range("I21").Value = "Nee"
$logs.Range("J21").Value = "Nee"

# --- Logs sheet: extend conditional formatting ranges down to row 21 -----
$logs.Range("D2:D20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D21"))
$logs.Range("G2:G20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G21"))
$logs.Range("H2:H20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H21"))
$logs.Range("I2:I20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I21"))
$logs.Range("J2:J20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J21"))

# --- Dashboard sheet: recount categories, swapping rows 3 and 4 ----------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Intern verzoek / Actie voor medewerker"
$dash.Range("B3").Value = 5
$dash.Range("A4").Value = "Productinformatie"
$dash.Range("B4").Value = 5
